$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Insert a new row above row 30, shifting existing rows (30-43) down to (31-44).
$ws.Rows.Item(30).Insert()

# Fill in the newly inserted row 30 with the "Meeting with Amit" entry (2018-04-11, 1 hour).
$ws.Cells.Item(30, 1).Value = 43201
$ws.Cells.Item(30, 2).Value = "Meeting with Amit"
$ws.Cells.Item(30, 3).Value = 1

# Fill in row 36 (previously empty, now holding the last real entry after the shift)
# with the new "github" entry (2018-05-02, 4 hours).
$ws.Cells.Item(36, 1).Value = 43222
$ws.Cells.Item(36, 2).Value = "Working on github and question 3 with Anna"
$ws.Cells.Item(36, 3).Value = 4

$ws.Range("M29").Select()
